# Append rows 38-68 (particella records 36-66) to Sheet1, matching the
# existing table's layout: col A = sequential index (bold/border/centered
# style copied from the existing A2 cell), col B = "codice_particella" as
# literal text, col C = "codice_comune_catastale" as a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    ,@(38, 36, "673/2", 384)
    ,@(39, 37, ".372", 384)
    ,@(40, 38, ".373", 384)
    ,@(41, 39, ".374", 384)
    ,@(42, 40, "673/2", 384)
    ,@(43, 41, "406/3", 384)
    ,@(44, 42, "605", 384)
    ,@(45, 43, "657/1", 384)
    ,@(46, 44, "674", 384)
    ,@(47, 45, "765/3", 384)
    ,@(48, 46, "938", 384)
    ,@(49, 47, "996", 384)
    ,@(50, 48, "2074", 384)
    ,@(51, 49, "2050", 384)
    ,@(52, 50, "2065", 384)
    ,@(53, 51, "2066", 384)
    ,@(54, 52, "2153", 384)
    ,@(55, 53, "2154", 384)
    ,@(56, 54, "1419", 287)
    ,@(57, 55, "1420", 287)
    ,@(58, 56, "1421", 287)
    ,@(59, 57, "1430", 287)
    ,@(60, 58, "1431/1", 287)
    ,@(61, 59, "1431/34", 287)
    ,@(62, 60, ".950/1", 287)
    ,@(63, 61, ".950/2", 287)
    ,@(64, 62, "239", 287)
    ,@(65, 63, "241", 287)
    ,@(66, 64, "773", 441)
    ,@(67, 65, ".4046", 307)
    ,@(68, 66, ".4047", 307)
)

# Grab the formatting of the first data row (bold/bordered/centered index
# column, plain text/number columns) so the new rows look identical to the
# existing ones.
$ws.Cells.Item(2, 1).Copy() | Out-Null

foreach ($r in $rows) {
    $rowNum = $r[0]
    $a = $r[1]
    $b = $r[2]
    $c = $r[3]

    # Column A: numeric index, same bold/border/centered style as the rest
    # of the table (copy-format so we reuse the existing style, not a new
    # one).
    $ws.Cells.Item($rowNum, 1).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($rowNum, 1).Value = $a

    # Column B: force text so values like ".372" / "605" aren't
    # reinterpreted as numbers, then restore the plain (unstyled) format so
    # no stray formatting is left behind.
    $ws.Cells.Item($rowNum, 2).NumberFormat = "@"
    $ws.Cells.Item($rowNum, 2).Value = $b
    $ws.Cells.Item($rowNum, 2).Style = $ws.Cells.Item(2, 2).Style

    # Column C: plain number.
    $ws.Cells.Item($rowNum, 3).Value = $c
}

$excel.CutCopyMode = $false
